# Rewrite the "Intro" paragraph of the ECB background article.
#
# Strategy: Word/this COM-interop runtime automatically coalesces two
# text-adjacent runs that end up with identical run properties (rPr)
# whenever Find.Execute mutates text touching either of them. Since the
# final paragraph has eight runs that all share the exact same rPr
# (rFonts "Civil Premium" + lang nl-NL), doing the text edits *after*
# unifying formatting would collapse everything into a single run.
#
# So: 1) do every text edit first, while runs still carry their original
#        (differing) formatting so natural run boundaries are preserved;
#        any run whose neighbour would coalesce away a boundary we still
#        need gets a throwaway distinguishing font first;
#     2) only at the very end, once all text + boundaries are correct,
#        reformat the whole paragraph range (including the paragraph
#        mark) in one shot -- a pure formatting assignment does not
#        trigger the coalescing behaviour, so the eight runs survive.

$d = $word.ActiveDocument

# Locate the paragraph: "Intro: in de laatste 3 jaar hebben we een enorme
# verandering  Op deze veranderingen ... ###Betere BU hier"
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "Intro:*") {
        $target = $para
        break
    }
}
$p = $target

# --- 1. Text edits -------------------------------------------------

# "Intro: " (run 1) -> "D" ; give it a throwaway font first so the
# replace doesn't get merged into the still-unedited run 2 that follows
# (both currently share the same rPr apart from rsid bookkeeping). Each
# Find.Execute needs its own fresh range -- re-using a range whose
# current extent already equals the search hit makes Find search forward
# from the end of that range instead of matching in place.
$run1 = $p.Range.Duplicate
[void]$run1.Find.Execute("Intro: ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$run1.Font.Name = "Arial"
$run1b = $p.Range.Duplicate
[void]$run1b.Find.Execute("Intro: ", $false, $false, $false, $false, $false, $true, 1, $false, "D", 2)

# Old run 2 "in de laatste 3 jaar hebben we een enorme verandering " loses
# its "in d" prefix, keeping the rest in the same (still separate) run.
$e1 = $p.Range.Duplicate
[void]$e1.Find.Execute("in d", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

# ... then its tail is rewritten in place (still the same run as above).
$e2 = $p.Range.Duplicate
[void]$e2.Find.Execute("hebben we een enorme verandering ", $false, $false, $false, $false, $false, $true, 1, $false, "zijn er economisch veel harde klappen gevoeld. ", 2)

# Insert a brand new sentence right before the old single-space run
# (run 3, still highlighted yellow at this point) so it lands as its own
# run rather than merging into the sentence we just edited above.
$insPoint = $p.Range.Duplicate
[void]$insPoint.Find.Execute(" Op deze veranderingen", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insPoint.Collapse(1)
$insPoint.InsertBefore("De Nederlandse burger weet eigenlijk helemaal niet wie er nou verantwoordelijk is voor de stijging van hun boodschappenmandje.")

# Give the freshly inserted sentence a throwaway font so subsequent edits
# can't accidentally fold it back into its neighbours.
$newSentence = $p.Range.Duplicate
[void]$newSentence.Find.Execute("De Nederlandse burger weet eigenlijk helemaal niet wie er nou verantwoordelijk is voor de stijging van hun boodschappenmandje.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$newSentence.Font.Name = "Arial Black"

# Old run 4 ("Op deze veranderingen ... ###Betere BU hier") is rewritten
# into three pieces: a rephrased prefix, the untouched "Europese Centrale
# Bank" phrase, and a rephrased suffix. Mark each piece with its own
# throwaway font first so the three final runs stay distinct from one
# another (and from the still-highlighted space run before them).
$prefix = $p.Range.Duplicate
[void]$prefix.Find.Execute("Op deze veranderingen van het Nederlandse prijspeil heeft de ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$prefix.Font.Name = "Courier New"

$mid = $p.Range.Duplicate
[void]$mid.Find.Execute("Europese Centrale Bank", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$mid.Font.Name = "Times New Roman"

$suffix = $p.Range.Duplicate
[void]$suffix.Find.Execute(" (ECB) een gigantische invloed. ###Betere BU hier", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$suffix.Font.Name = "Georgia"

# Now rewrite the prefix and suffix text in place (mid is left untouched).
$prefix2 = $p.Range.Duplicate
[void]$prefix2.Find.Execute("Op deze veranderingen van het Nederlandse prijspeil heeft de ", $false, $false, $false, $false, $false, $true, 1, $false, "De keuzes die de ", 2)

$suffix2 = $p.Range.Duplicate
[void]$suffix2.Find.Execute(" (ECB) een gigantische invloed. ###Betere BU hier", $false, $false, $false, $false, $false, $true, 1, $false, " ECB laat gemaakt heeft hebben een gigantische invloed hierop gehad.", 2)

# --- 2. Final formatting pass ---------------------------------------

# Unify every run (and the paragraph mark) onto "Civil Premium" / nl-NL,
# and strip the leftover yellow highlight. A plain formatting assignment
# (as opposed to a Find.Execute replace) does not coalesce runs, so the
# eight text runs built above remain intact.
$finalRange = $p.Range.Duplicate
$finalRange.Font.Name = "Civil Premium"
$finalRange.HighlightColorIndex = 0

Write-Output $p.Range.Text
